# Apply "밀린 커밋" changes: mark progress ("O") for rows 16 & 17
# Row 16 ("21. 자바 ORM 표준 JPA 프로그래밍 - 기본편"): mark C16:E16 as done,
#   and switch B16 / O16 to the "no data" (X) style used elsewhere on the sheet.
# Row 17 ("22. 실전! 스프링 부트와 JPA 활용1 - 웹 어플리케이션 개발"): mark
#   F17:I17 as done (matching the already-marked B17:E17) and restyle the
#   row label A17 to the highlighted style used by fully-marked rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Row 16 -------------------------------------------------------------

# B16 & O16 become the bordered "X" style (same style already used by P16)
$ws.Range("P16").Copy()
$ws.Range("B16").PasteSpecial($xlPasteFormats)
$ws.Range("O16").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# C16:E16 become highlighted "done" cells containing "O" (copy style from
# an existing highlighted cell, e.g. B7, then set the value)
$ws.Range("B7").Copy()
$ws.Range("C16:E16").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("C16").Value = "O"
$ws.Range("D16").Value = "O"
$ws.Range("E16").Value = "O"

# --- Row 17 -------------------------------------------------------------

# A17 label gets the highlighted row-label style (same as A7/A8)
$ws.Range("A7").Copy()
$ws.Range("A17").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# F17:I17 become highlighted "done" cells containing "O", matching B17:E17
$ws.Range("B17").Copy()
$ws.Range("F17:I17").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("F17").Value = "O"
$ws.Range("G17").Value = "O"
$ws.Range("H17").Value = "O"
$ws.Range("I17").Value = "O"
